# Applies the scheduled "Updated cryptos list" GitHub Actions refresh to
# cryptos.xlsx: refreshed Price (col D) / Volume(1h) (col E) text for most
# rows, plus a few adjacent row pairs that swapped rank (Coin/Link/Price/
# Volume all move together).
#
# Every written cell is literal text in the source (t="inlineStr"), including
# column D where the text often looks like a plain number (e.g. "1.00",
# "52.630.01"). Excel's Range.Value auto-coerces such strings to real
# numbers, so for column D we briefly force a text NumberFormat while
# assigning, then restore "General" so the cell style matches the original
# (only the stored type/value changes, not formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.630.01"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -12.81%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.338.61"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -19.04%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.05%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "439.60"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -16.05%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.29"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -12.76%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.48%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.472"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -13.87%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.345.48"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -18.91%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.26"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -12.30%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0898"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -16.23%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.308"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -13.78%  "
# Row 13
$ws.Range("E13").Value = "  -5.32%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "52.709.25"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -12.74%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.17"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -15.35%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000120"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -14.19%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.353.12"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -18.96%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.03"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -18.70%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "304.04"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -14.40%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.10"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -21.55%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.01%  "
# Row 22
$ws.Range("E22").Value = "  -1.91%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.24"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -20.69%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "54.13"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -16.15%  "
# Row 25
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.151"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -15.66%  "
# Row 26
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.371"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -17.78%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.18"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -7.96%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.47%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0687"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -17.70%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "144.47"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.67%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "17.20"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -12.23%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.34"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -19.72%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.83"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -12.91%  "
# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.840"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -15.50%  "
# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.56"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -17.38%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -15.05%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.996"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.22%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "32.23"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -14.24%  "
# Row 39
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.20"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -13.37%  "
# Row 40
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.14"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -1.58%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0509"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -12.32%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.24"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -16.09%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.932.85"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -15.31%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.529"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -17.92%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0210"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -10.89%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0837"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -8.62%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.10"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -16.19%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.91"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -20.69%  "
# Row 49
$ws.Range("E49").Value = "  -5.75%  "
# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.53"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -12.06%  "
# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "15.33"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -15.47%  "
